# actualizacao da sessao 8 chimanimani
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    # Force the cell to keep a literal text value (shared string) instead of
    # letting Excel auto-convert numeric- or date-looking strings into a
    # number / date serial. We briefly apply a text number format, assign
    # the value, then restore the cell style back to Normal/General so the
    # saved style index matches the original (unformatted) cells.
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2 - AFONSO BINGALA JONE
Set-TextValue "D2" "GORONGOSA"
Set-TextValue "E2" "TECNICO"
Set-TextValue "F2" "845123123"
Set-TextValue "G2" "10/22/2022"
$ws.Range("H2").Value = 4
Set-TextValue "K2" "863033655"

# Row 3 - ANTONIO AGOSTINHO JOAO NOBRE
Set-TextValue "D3" "GORONGOSA"
Set-TextValue "E3" "TECNICO DO CAMPO"
Set-TextValue "F3" "848226339"
Set-TextValue "G3" "10/21/2022"
$ws.Range("H3").Value = 3
Set-TextValue "K3" "841589157"

# Row 5 - JOSSEFO CELESTINO SALIVA
Set-TextValue "D5" "PARQUE DE GORONGOSA"
Set-TextValue "E5" "TECNICO DE CAFE"
$ws.Range("F5").Value = ""
Set-TextValue "G5" " 11/9/2022"
$ws.Range("H5").Value = 2

# Row 6 - JULEIDA ZULFA CARLOS
Set-TextValue "G6" " 11/9/2022"
$ws.Range("H6").Value = 2
